# Add two new income rows (Salary 20000 on 2024-05-07 w/ timestamp, and a
# text-dated Salary 90 row) to the "Rafid" income-record sheet, matching the
# upstream "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing date cells (A2:A5) carry a custom datetime number format whose
# formatCode is effectively duplicated (numFmtId 164 lowercase / 165
# uppercase) in the original workbook. Re-applying the lowercase form here
# is what the authoring tool did when it touched this column, and it is
# also what lets the freshly-added A6 date cell below pick up the very same
# style as the pre-existing rows.
$ws.Range("A2:A6").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Row 6: a new Salary entry, timestamped (fractional serial date).
$ws.Range("A6").Value = 45419.96802559028
$ws.Range("B6").Value = "Salary"
$ws.Range("C6").Value = 20000

# Row 7: another Salary entry, but the date column here was uploaded as a
# literal text string ("2024-05-07"), not a real date serial. Force the
# cell to text first so Excel doesn't auto-convert the string into a date
# value, then strip the formatting back off so the cell ends up with no
# explicit style, just like the source file.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2024-05-07"
$ws.Range("A7").ClearFormats()
$ws.Range("B7").Value = "Salary"
$ws.Range("C7").Value = 90
